$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A82").Value = "teste"
$ws.Range("B82").Value = "Incompleto"
$ws.Range("C82").Value = "PS3"
$ws.Range("D82").Value = "Zerar"
